$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 8 (the formula row referencing A3); this shifts rows 9-14 up to 8-13
$ws.Rows.Item(8).Delete()
